$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 4056.1667
$ws.Range("I11").Value = 4056.1667
$ws.Range("K11").Value = 4056.1667
$ws.Range("M11").Value = -3916.1667

# Row 53
$ws.Range("H53").Value = 837.4
$ws.Range("I53").Value = 905
$ws.Range("J53").Value = 769.8
$ws.Range("K53").Value = 905
$ws.Range("L53").Value = 769.8
$ws.Range("M53").Value = -268
$ws.Range("N53").Value = -2043.8

# Row 70
$ws.Range("H70").Value = 38892224
$ws.Range("I70").Value = 41667420
$ws.Range("J70").Value = 37042096
$ws.Range("K70").Value = 125002260
$ws.Range("L70").Value = 111126288
$ws.Range("M70").Value = -125001990
$ws.Range("N70").Value = -111126828

# Row 73
$ws.Range("H73").Value = 38892224
$ws.Range("I73").Value = 41667420
$ws.Range("J73").Value = 37042096
$ws.Range("K73").Value = 125002260
$ws.Range("L73").Value = 111126288
$ws.Range("M73").Value = -125001324
$ws.Range("N73").Value = -111128160

# Row 106
$ws.Range("H106").Value = 2151
$ws.Range("I106").Value = 1974.6154
$ws.Range("K106").Value = 1974.6154
$ws.Range("M106").Value = -1343.6154

# Row 107
$ws.Range("H107").Value = 18480750
$ws.Range("J107").Value = 18183662
$ws.Range("L107").Value = 18183662
$ws.Range("N107").Value = -18187502

# Row 110
$ws.Range("H110").Value = 91990
$ws.Range("J110").Value = 91990
$ws.Range("L110").Value = 91990
$ws.Range("N110").Value = -100170

# Row 112
$ws.Range("H112").Value = 13454.5
$ws.Range("J112").Value = 13454.5
$ws.Range("L112").Value = 40363.5
$ws.Range("N112").Value = -42579.5

# Row 137
$ws.Range("H137").Value = 3043.1316
$ws.Range("I137").Value = 3032.3635
$ws.Range("J137").Value = 3057.9375
$ws.Range("K137").Value = 9097.0905
$ws.Range("L137").Value = 9173.8125
$ws.Range("M137").Value = -6547.0905
$ws.Range("N137").Value = -14273.8125

# Row 138
$ws.Range("H138").Value = 1891487.5
$ws.Range("I138").Value = 4551.75
$ws.Range("J138").Value = 2443761.2
$ws.Range("K138").Value = 13655.25
$ws.Range("L138").Value = 7331283.600000001
$ws.Range("M138").Value = -8515.25
$ws.Range("N138").Value = -7341563.600000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1789623.9
$ws.Range("I32").Value = 2196442.2
$ws.Range("K32").Value = 2196442.2
$ws.Range("M32").Value = -2196155.2

# Row 61
$ws.Range("H61").Value = 6369.6665
$ws.Range("I61").Value = 3475.5
$ws.Range("K61").Value = 3475.5
$ws.Range("M61").Value = -3263.5

# Row 132
$ws.Range("H132").Value = 2033869
$ws.Range("I132").Value = 5858985
$ws.Range("J132").Value = 8807.529
$ws.Range("K132").Value = 17576955
$ws.Range("L132").Value = 26422.587
$ws.Range("M132").Value = -17574425
$ws.Range("N132").Value = -31482.587

# Row 136
$ws.Range("H136").Value = 6369.6665
$ws.Range("I136").Value = 3475.5
$ws.Range("K136").Value = 10426.5
$ws.Range("M136").Value = -7876.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3933.55
$ws.Range("I105").Value = 2731.7273
$ws.Range("J105").Value = 5402.4443
$ws.Range("K105").Value = 2731.7273
$ws.Range("L105").Value = 5402.4443
$ws.Range("M105").Value = -984.7273
$ws.Range("N105").Value = -8896.444299999999

# Row 134
$ws.Range("H134").Value = 6159.613
$ws.Range("I134").Value = 1095.4615
$ws.Range("K134").Value = 3286.3845
$ws.Range("M134").Value = -751.3844999999997

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4447.971
$ws.Range("I16").Value = 3872.0386
$ws.Range("K16").Value = 3872.0386
$ws.Range("M16").Value = -3585.0386

# Row 19
$ws.Range("H19").Value = 976
$ws.Range("I19").Value = 976
$ws.Range("K19").Value = 976
$ws.Range("M19").Value = -806

# Row 24
$ws.Range("H24").Value = 976
$ws.Range("I24").Value = 976
$ws.Range("K24").Value = 976
$ws.Range("M24").Value = -806

# Row 31
$ws.Range("H31").Value = 8837.5
$ws.Range("I31").Value = 5117.067
$ws.Range("K31").Value = 5117.067
$ws.Range("M31").Value = -4822.067

# Row 34
$ws.Range("H34").Value = 8837.5
$ws.Range("I34").Value = 5117.067
$ws.Range("K34").Value = 5117.067
$ws.Range("M34").Value = -4915.067

# Row 113
$ws.Range("H113").Value = 4447.971
$ws.Range("I113").Value = 3872.0386
$ws.Range("K113").Value = 3872.0386
$ws.Range("M113").Value = -1702.0386

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 87730.914
$ws.Range("J2").Value = 287367.72
$ws.Range("L2").Value = 1724206.32
$ws.Range("N2").Value = -1724432.32

# Row 3
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 6000
$ws.Range("M3").Value = -5888

# Row 5
$ws.Range("H5").Value = 1335912.1
$ws.Range("J5").Value = 4274.1113
$ws.Range("L5").Value = 12822.3339
$ws.Range("N5").Value = -13046.3339

# Row 68
$ws.Range("H68").Value = 2892.8164
$ws.Range("J68").Value = 3017.8333
$ws.Range("L68").Value = 9053.499899999999
$ws.Range("N68").Value = -10675.4999

# Row 71
$ws.Range("H71").Value = 2892.8164
$ws.Range("J71").Value = 3017.8333
$ws.Range("L71").Value = 27160.4997
$ws.Range("N71").Value = -35272.4997

# Row 92
$ws.Range("H92").Value = 5918113.5
$ws.Range("I92").Value = 1144
$ws.Range("J92").Value = 7693204.5
$ws.Range("K92").Value = 3432
$ws.Range("L92").Value = 23079613.5
$ws.Range("M92").Value = -2184
$ws.Range("N92").Value = -23082109.5

# Row 135
$ws.Range("H135").Value = 1335912.1
$ws.Range("J135").Value = 4274.1113
$ws.Range("L135").Value = 38467.00169999999
$ws.Range("N135").Value = -43537.00169999999

# Row 137
$ws.Range("H137").Value = 141512.8
$ws.Range("I137").Value = 102069.4
$ws.Range("J137").Value = 220399.6
$ws.Range("K137").Value = 306208.2
$ws.Range("L137").Value = 661198.8
$ws.Range("M137").Value = -301108.2
$ws.Range("N137").Value = -671398.8

# Row 139
$ws.Range("H139").Value = 279318
$ws.Range("I139").Value = 1002499.7
$ws.Range("K139").Value = 3007499.1
$ws.Range("M139").Value = -3002359.1

# Row 140
$ws.Range("H140").Value = 224318.94
$ws.Range("I140").Value = 309141.7
$ws.Range("K140").Value = 927425.1000000001
$ws.Range("M140").Value = -922245.1000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 12997.5
$ws.Range("I132").Value = 1989
$ws.Range("K132").Value = 5967
$ws.Range("M132").Value = -3437

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7916.3
$ws.Range("I7").Value = 6892.1665
$ws.Range("J7").Value = 8355.214
$ws.Range("K7").Value = 6892.1665
$ws.Range("L7").Value = 8355.214
$ws.Range("M7").Value = -6780.1665
$ws.Range("N7").Value = -8579.214

# Row 46
$ws.Range("H46").Value = 3779.9614
$ws.Range("I46").Value = 4451.4
$ws.Range("J46").Value = 3620.0952
$ws.Range("K46").Value = 4451.4
$ws.Range("L46").Value = 3620.0952
$ws.Range("M46").Value = -4263.4
$ws.Range("N46").Value = -3996.0952

# Row 100
$ws.Range("H100").Value = 4483.8335
$ws.Range("I100").Value = 3979.8
$ws.Range("J100").Value = 7004
$ws.Range("K100").Value = 3979.8
$ws.Range("L100").Value = 7004
$ws.Range("M100").Value = -3438.8
$ws.Range("N100").Value = -8086

# Row 126
$ws.Range("H126").Value = 7916.3
$ws.Range("I126").Value = 6892.1665
$ws.Range("J126").Value = 8355.214
$ws.Range("K126").Value = 20676.4995
$ws.Range("L126").Value = 25065.642
$ws.Range("M126").Value = -18206.4995
$ws.Range("N126").Value = -30005.642

# Row 127
$ws.Range("H127").Value = 59378.5
$ws.Range("J127").Value = 59378.5
$ws.Range("L127").Value = 59378.5
$ws.Range("N127").Value = -69298.5

# Row 132
$ws.Range("H132").Value = 16136081
$ws.Range("I132").Value = 31253352
$ws.Range("K132").Value = 93760056
$ws.Range("M132").Value = -93757526

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 5712.2354
$ws.Range("J122").Value = 5614.4
$ws.Range("L122").Value = 16843.2
$ws.Range("N122").Value = -21743.2

# Row 125
$ws.Range("H125").Value = 50913
$ws.Range("J125").Value = 50913
$ws.Range("L125").Value = 50913
$ws.Range("N125").Value = -60753

# Row 132
$ws.Range("H132").Value = 27796812
$ws.Range("I132").Value = 45459156
$ws.Range("K132").Value = 136377468
$ws.Range("M132").Value = -136374938

# Row 133
$ws.Range("H133").Value = 144975
$ws.Range("J133").Value = 144975
$ws.Range("L133").Value = 144975
$ws.Range("N133").Value = -155095

# Row 136
$ws.Range("H136").Value = 25644440
$ws.Range("I136").Value = 41668456
$ws.Range("J136").Value = 6016.8667
$ws.Range("K136").Value = 125005368
$ws.Range("L136").Value = 18050.6001
$ws.Range("M136").Value = -125002818
$ws.Range("N136").Value = -23150.6001
